$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Date column (BF) keeps its text format so values are not
# auto-converted into date serial numbers.
$ws.Range("BF2:BF31").NumberFormat = "@"

$ws.Range("AD2").Value = 20
$ws.Range("AR2").Value = 26
$ws.Range("BF2").Value = "2014-02-06"
$ws.Range("AQ3").Value = 12
$ws.Range("BF3").Value = "2014-02-06"
$ws.Range("D4").Value = 46
$ws.Range("E4").Value = 21
$ws.Range("G4").Value = 0.457
$ws.Range("I4").Value = 35.2
$ws.Range("J4").Value = 78.2
$ws.Range("K4").Value = 0.45
$ws.Range("L4").Value = 8
$ws.Range("M4").Value = 21.7
$ws.Range("O4").Value = 19.2
$ws.Range("P4").Value = 24.9
$ws.Range("Q4").Value = 0.771
$ws.Range("R4").Value = 9.1
$ws.Range("S4").Value = 29.8
$ws.Range("T4").Value = 38.9
$ws.Range("U4").Value = 20.7
$ws.Range("Z4").Value = 22.3
$ws.Range("AA4").Value = 21.3
$ws.Range("AB4").Value = 97.5
$ws.Range("AC4").Value = -2.7
$ws.Range("AD4").Value = 30
$ws.Range("AE4").Value = 19
$ws.Range("AF4").Value = 17
$ws.Range("AI4").Value = 28
$ws.Range("AM4").Value = 14
$ws.Range("AO4").Value = 6
$ws.Range("AR4").Value = 27
$ws.Range("AS4").Value = 28
$ws.Range("AU4").Value = 19
$ws.Range("AY4").Value = 7
$ws.Range("BA4").Value = 10
$ws.Range("BC4").Value = 22
$ws.Range("BF4").Value = "2014-02-06"
$ws.Range("AI5").Value = 27
$ws.Range("BC5").Value = 21
$ws.Range("BF5").Value = "2014-02-06"
$ws.Range("D6").Value = 48
$ws.Range("F6").Value = 24
$ws.Range("G6").Value = 0.5
$ws.Range("L6").Value = 5.9
$ws.Range("M6").Value = 17.5
$ws.Range("P6").Value = 23.6
$ws.Range("Q6").Value = 0.765
$ws.Range("R6").Value = 12
$ws.Range("V6").Value = 15.6
$ws.Range("X6").Value = 5.3
$ws.Range("Z6").Value = 19.4
$ws.Range("AA6").Value = 21.4
$ws.Range("AB6").Value = 92.2
$ws.Range("AC6").Value = -0.5
$ws.Range("AD6").Value = 20
$ws.Range("AF6").Value = 15
$ws.Range("AG6").Value = 15
$ws.Range("AJ6").Value = 26
$ws.Range("AQ6").Value = 13
$ws.Range("BA6").Value = 7
$ws.Range("BF6").Value = "2014-02-06"
$ws.Range("AD7").Value = 9
$ws.Range("AH7").Value = 3
$ws.Range("BF7").Value = "2014-02-06"
$ws.Range("AE8").Value = 8
$ws.Range("BF8").Value = "2014-02-06"
$ws.Range("AO9").Value = 7
$ws.Range("BA9").Value = 8
$ws.Range("BB9").Value = 10
$ws.Range("BF9").Value = "2014-02-06"
$ws.Range("AD10").Value = 20
$ws.Range("BF10").Value = "2014-02-06"
$ws.Range("D11").Value = 49
$ws.Range("E11").Value = 29
$ws.Range("G11").Value = 0.592
$ws.Range("J11").Value = 85
$ws.Range("M11").Value = 24.4
$ws.Range("N11").Value = 0.38
$ws.Range("O11").Value = 16.1
$ws.Range("P11").Value = 21.8
$ws.Range("S11").Value = 34.9
$ws.Range("T11").Value = 46.1
$ws.Range("U11").Value = 22.8
$ws.Range("V11").Value = 16.3
$ws.Range("Z11").Value = 22.1
$ws.Range("AA11").Value = 19.9
$ws.Range("AB11").Value = 103
$ws.Range("AC11").Value = 4.1
$ws.Range("AD11").Value = 9
$ws.Range("AH11").Value = 17
$ws.Range("AK11").Value = 10
$ws.Range("AP11").Value = 20
$ws.Range("BC11").Value = 8
$ws.Range("BF11").Value = "2014-02-06"
$ws.Range("AH12").Value = 20
$ws.Range("BF12").Value = "2014-02-06"
$ws.Range("AD13").Value = 20
$ws.Range("AJ13").Value = 27
$ws.Range("BF13").Value = "2014-02-06"
$ws.Range("BF14").Value = "2014-02-06"
$ws.Range("AD15").Value = 9
$ws.Range("AT15").Value = 21
$ws.Range("AU15").Value = 9
$ws.Range("AY15").Value = 11
$ws.Range("BF15").Value = "2014-02-06"
$ws.Range("AD16").Value = 20
$ws.Range("AK16").Value = 9
$ws.Range("BF16").Value = "2014-02-06"
$ws.Range("AD17").Value = 20
$ws.Range("AG17").Value = 4
$ws.Range("AL17").Value = 13
$ws.Range("BF17").Value = "2014-02-06"
$ws.Range("AD18").Value = 9
$ws.Range("AO18").Value = 28
$ws.Range("AU18").Value = 18
$ws.Range("BF18").Value = "2014-02-06"
$ws.Range("AD19").Value = 9
$ws.Range("AF19").Value = 17
$ws.Range("AG19").Value = 17
$ws.Range("BC19").Value = 7
$ws.Range("BF19").Value = "2014-02-06"
$ws.Range("AD20").Value = 20
$ws.Range("AE20").Value = 19
$ws.Range("BF20").Value = "2014-02-06"
$ws.Range("AD21").Value = 9
$ws.Range("AS21").Value = 27
$ws.Range("BF21").Value = "2014-02-06"
$ws.Range("BF22").Value = "2014-02-06"
$ws.Range("AJ23").Value = 21
$ws.Range("AR23").Value = 25
$ws.Range("BF23").Value = "2014-02-06"
$ws.Range("AM24").Value = 15
$ws.Range("BF24").Value = "2014-02-06"
$ws.Range("AD25").Value = 9
$ws.Range("AE25").Value = 8
$ws.Range("AG25").Value = 8
$ws.Range("AK25").Value = 8
$ws.Range("BA25").Value = 11
$ws.Range("BF25").Value = "2014-02-06"
$ws.Range("AD26").Value = 9
$ws.Range("AF26").Value = 5
$ws.Range("BF26").Value = "2014-02-06"
$ws.Range("AD27").Value = 9
$ws.Range("BF27").Value = "2014-02-06"
$ws.Range("D28").Value = 49
$ws.Range("F28").Value = 13
$ws.Range("G28").Value = 0.735
$ws.Range("I28").Value = 40.5
$ws.Range("J28").Value = 82.40000000000001
$ws.Range("K28").Value = 0.491
$ws.Range("N28").Value = 0.394
$ws.Range("O28").Value = 15.1
$ws.Range("P28").Value = 19.7
$ws.Range("Q28").Value = 0.768
$ws.Range("S28").Value = 33.5
$ws.Range("T28").Value = 42.5
$ws.Range("U28").Value = 24.9
$ws.Range("X28").Value = 5
$ws.Range("Y28").Value = 4.9
$ws.Range("Z28").Value = 18.1
$ws.Range("AA28").Value = 19.4
$ws.Range("AB28").Value = 104.2
$ws.Range("AC28").Value = 7
$ws.Range("AD28").Value = 9
$ws.Range("AF28").Value = 3
$ws.Range("AG28").Value = 3
$ws.Range("AH28").Value = 17
$ws.Range("AJ28").Value = 20
$ws.Range("AL28").Value = 12
$ws.Range("AO28").Value = 29
$ws.Range("AR28").Value = 28
$ws.Range("AT28").Value = 20
$ws.Range("AX28").Value = 11
$ws.Range("BB28").Value = 9
$ws.Range("BF28").Value = "2014-02-06"
$ws.Range("AD29").Value = 9
$ws.Range("AU29").Value = 17
$ws.Range("BF29").Value = "2014-02-06"
$ws.Range("AD30").Value = 20
$ws.Range("AP30").Value = 21
$ws.Range("BF30").Value = "2014-02-06"
$ws.Range("AD31").Value = 20
$ws.Range("AY31").Value = 8
$ws.Range("BF31").Value = "2014-02-06"
